$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds plain numeric-looking text (e.g. "3.400") where trailing
# zeros / exact formatting matter, so force text format before assigning
# to avoid Excel auto-converting the input into a numeric value.
$priceCells = @("D2","D3","D4","D5","D6","D7","D8","D9","D10","D11","D12","D13","D14","D15","D16","D17","D18","D19","D20","D21","D22","D23","D24","D25","D40","D41","D43","D44","D45","D47","D48","D50","D51")
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# --- Price (column D) updates ---
$ws.Range("D2").Value = "242.27"
$ws.Range("D3").Value = "22.88"
$ws.Range("D4").Value = "5.373"
$ws.Range("D5").Value = "0.05945"
$ws.Range("D6").Value = "3.401"
$ws.Range("D7").Value = "6.480"
$ws.Range("D8").Value = "0.8049"
$ws.Range("D9").Value = "0.9062"
$ws.Range("D10").Value = "0.1418"
$ws.Range("D11").Value = "0.07408"
$ws.Range("D12").Value = "0.03309"
$ws.Range("D13").Value = "0.03032"
$ws.Range("D14").Value = "0.09323"
$ws.Range("D15").Value = "3.849"
$ws.Range("D16").Value = "0.001582"
$ws.Range("D17").Value = "0.04519"
$ws.Range("D18").Value = "0.006098"
$ws.Range("D19").Value = "0.005000"
$ws.Range("D20").Value = "0.007494"
$ws.Range("D21").Value = "0.0009859"
$ws.Range("D22").Value = "0.00007801"
$ws.Range("D23").Value = "3.611"
$ws.Range("D24").Value = "2.136"
$ws.Range("D25").Value = "0.01115"
$ws.Range("D40").Value = "0.03852"
$ws.Range("D41").Value = "0.006068"
$ws.Range("D43").Value = "0.002800"
$ws.Range("D44").Value = "0.007169"
$ws.Range("D45").Value = "0.00005183"
$ws.Range("D47").Value = "0.0005801"
$ws.Range("D48").Value = "0.9701"
$ws.Range("D50").Value = "0.00002100"
$ws.Range("D51").Value = "0.0002000"

# Restore the default General format now that the text values are set,
# so the cells keep rendering/sorting like the rest of the sheet.
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "General"
}

# --- Coin / Link / Volume(1h) label updates (rows 18-25 coin re-ranking) ---
$ws.Range("B18").Value = "TigerCash"
$ws.Range("C18").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("E18").Value = "17TigerCashTCH"
$ws.Range("B19").Value = "HotbitToken"
$ws.Range("C19").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
$ws.Range("E19").Value = "18HotbitTokenHTB"
$ws.Range("B20").Value = "UpBots"
$ws.Range("C20").Value = "https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt"
$ws.Range("E20").Value = "19UpBotsUBXTBestin24h"
$ws.Range("B21").Value = "BitKan"
$ws.Range("C21").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
$ws.Range("E21").Value = "20BitKanKAN"
$ws.Range("B22").Value = "NitroEx"
$ws.Range("C22").Value = "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
$ws.Range("E22").Value = "21NitroExNTX"
$ws.Range("B23").Value = "LEO"
$ws.Range("C23").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("E23").Value = "22LEOLEO"
$ws.Range("B24").Value = "BTSEToken"
$ws.Range("C24").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("E24").Value = "23BTSETokenBTSE"
$ws.Range("B25").Value = "One"
$ws.Range("C25").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("E25").Value = "24OneONE"

Write-Host "Applied cryptos.xlsx symbol-list update"
